$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15, shifting old rows 15-18 down to 16-19.
$ws.Rows.Item(15).Insert()

# Row 14 becomes the new "NR_hdp_gb_1" result row (SBS_set2 / Realistic).
$ws.Range("B14").Value2 = "NR_hdp_gb_1"
$ws.Range("C14").Value2 = 2.19707955583194
$ws.Range("D14").Value2 = 0.0773183260467158
$ws.Range("E14").Value2 = 0.319131092406954
$ws.Range("F14").Value2 = 0.89375
$ws.Range("G14").Value2 = 0.984198463424988

# Row 15 (newly inserted) holds what used to be row 14's "signeR" data.
# Restore the merged "SBS_set2" label text that the insert left blank.
$ws.Range("A15").Value2 = "SBS_set2"
$ws.Range("B15").Value2 = "signeR"
$ws.Range("C15").Value2 = 2.07081907525817
$ws.Range("D15").Value2 = 0.0824169433012825
$ws.Range("E15").Value2 = 0.57701888915682
$ws.Range("F15").Value2 = 0.5125
$ws.Range("G15").Value2 = 0.981300186101349

Write-Host "A14:" $ws.Range("A14").Value2
Write-Host "B14:" $ws.Range("B14").Value2
Write-Host "A15:" $ws.Range("A15").Value2
Write-Host "B15:" $ws.Range("B15").Value2
Write-Host "A16:" $ws.Range("A16").Value2
Write-Host "B16:" $ws.Range("B16").Value2
Write-Host "A17:" $ws.Range("A17").Value2
Write-Host "A18:" $ws.Range("A18").Value2
Write-Host "A19:" $ws.Range("A19").Value2
